# Commit: "properties to properties tabs of tourney sheets"
#
# The Tournament sheet's "tournament" table mixed true tournament facts
# (name/location/timezone/venue labels) together with lookup-key rows
# (competition-key, host-key, venue-key.1..8) that really belong with the
# Colors sheet's key/value data. This script pulls all of those key/value
# rows out into a brand-new "Properties" table/sheet (key, value, notes)
# and removes the now-redundant rows from the Tournament table.

$wb = $excel.ActiveWorkbook

$wsTournament = $wb.Worksheets.Item("Tournament")
$wsColors     = $wb.Worksheets.Item("Colors")

# ---------------------------------------------------------------------
# 1. Create the new "Properties" sheet (placed after "Colors") and its
#    key/value/notes table, gathering data that used to live scattered
#    across the Tournament ("*-key" rows) and Colors ("notes" column)
#    tables.
# ---------------------------------------------------------------------
$wsProperties = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsColors)
$wsProperties.Name = "Properties"

$propertyRows = @(
    @("key",         "value",             $null),
    @("competition",  "mens-world-cup",   $null),
    @("host",         "qatar",            $null),
    @("timezone",     "Asia/Qatar",       $null),
    @("color.a",      "#94d9f5",          "cyan"),
    @("color.b",      "#fee289",          "yellow"),
    @("color.c",      "#f79d8f",          "red"),
    @("color.d",      "#c4e1b5",          "green"),
    @("color.e",      "#b0d0ee",          "blue"),
    @("color.f",      "#c0e4df",          "teal"),
    @("color.g",      "#fab077",          "orange"),
    @("color.h",      "#eecbef",          "purple"),
    @("venue.01",     "qa-al-bayt",       $null),
    @("venue.02",     "qa-khalifa",       $null),
    @("venue.03",     "qa-al-thumama",    $null),
    @("venue.04",     "qa-ahmad-bin-ali", $null),
    @("venue.05",     "qa-lusail",        $null),
    @("venue.06",     "qa-974",           $null),
    @("venue.07",     "qa-education-city",$null),
    @("venue.08",     "qa-al-janoub",     $null)
)

for ($i = 0; $i -lt $propertyRows.Count; $i++) {
    $r = $i + 1
    $data = $propertyRows[$i]
    $wsProperties.Cells.Item($r, 1).Value = $data[0]
    $wsProperties.Cells.Item($r, 2).Value = $data[1]
    if ($data[2] -ne $null) {
        $wsProperties.Cells.Item($r, 3).Value = $data[2]
    }
}

$propertiesTable = $wsProperties.ListObjects.Add(1, $wsProperties.Range("A1:C20"), [System.Reflection.Missing]::Value, 1)
$propertiesTable.Name = "Properties"

# match the authored column widths (key / value / notes)
$wsProperties.Columns.Item(1).ColumnWidth = 9.830729166666666
$wsProperties.Columns.Item(2).ColumnWidth = 13.666666666666666
$wsProperties.Columns.Item(3).ColumnWidth = 7.166666666666667

$wsProperties.Range("A4:B4").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. Remove the rows from the Tournament table whose data now lives on
#    the Properties sheet: competition-key/host-key (rows 2-3) and the
#    venue-key.1..8 rows (rows 15-22). Delete bottom-up so row numbers
#    of not-yet-deleted rows stay stable.
# ---------------------------------------------------------------------
$tournamentRowsToDelete = @(22, 21, 20, 19, 18, 17, 16, 15, 3, 2)
foreach ($r in $tournamentRowsToDelete) {
    $wsTournament.Rows.Item($r).Delete()
}

$wsTournament.Activate()
$wsTournament.Range("A13:XFD20").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. Best-effort restore of the Colors sheet's last selection (a
#    non-contiguous highlight of its "value" and "notes" columns).
# ---------------------------------------------------------------------
$wsColors.Activate()
$wsColors.Range("J2:J9").Select() | Out-Null

# Leave the Tournament tab as the active one, same as before the edit.
$wsTournament.Activate()
